$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 for "Galba Freire Moita" (shifts rows 16-44 down to 17-45)
$ws.Rows.Item(16).Insert()

# Re-apply the full data range (Docente / status / Interno Fiocruz / Externo a Fiocruz)
# reflecting the October extraction update.
$ws.Cells.Item(2, 1).Value = "Adriana Costa Bacelo"
$ws.Cells.Item(2, 2).Value = "concluídas"
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(3, 1).Value = "Adriana Costa Bacelo"
$ws.Cells.Item(3, 2).Value = "em andamento"
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(4, 1).Value = "Ana Cláudia de Araújo Teixeira"
$ws.Cells.Item(4, 2).Value = "concluídas"
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(5, 1).Value = "Anya Pimentel Gomes Fernandes Vieira Meyer"
$ws.Cells.Item(5, 2).Value = "em andamento"
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(6, 1).Value = "Anya Pimentel Gomes Fernandes Vieira Meyer"
$ws.Cells.Item(6, 2).Value = "concluídas"
$ws.Cells.Item(6, 3).Value = 6
$ws.Cells.Item(6, 4).Value = 3
$ws.Cells.Item(7, 1).Value = "Carla Freire Celedonio Fernandes"
$ws.Cells.Item(7, 2).Value = "concluídas"
$ws.Cells.Item(7, 3).Value = 4
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(8, 1).Value = "Carla Freire Celedonio Fernandes"
$ws.Cells.Item(8, 2).Value = "em andamento"
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(9, 1).Value = "Donat Alexander de Chapeaurouge"
$ws.Cells.Item(9, 2).Value = "concluídas"
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(10, 1).Value = "Fabio Miyajima"
$ws.Cells.Item(10, 2).Value = "concluídas"
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(10, 4).Value = 4
$ws.Cells.Item(11, 1).Value = "Fabio Miyajima"
$ws.Cells.Item(11, 2).Value = "em andamento"
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(11, 4).Value = 7
$ws.Cells.Item(12, 1).Value = "Fernando Braga Stehling Dias"
$ws.Cells.Item(12, 2).Value = "concluídas"
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(13, 1).Value = "Fernando Braga Stehling Dias"
$ws.Cells.Item(13, 2).Value = "em andamento"
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 3
$ws.Cells.Item(14, 1).Value = "Fernando Ferreira Carneiro"
$ws.Cells.Item(14, 2).Value = "concluídas"
$ws.Cells.Item(14, 3).Value = 4
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(15, 1).Value = "Fernando Ferreira Carneiro"
$ws.Cells.Item(15, 2).Value = "em andamento"
$ws.Cells.Item(15, 3).Value = 3
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(16, 1).Value = "Galba Freire Moita"
$ws.Cells.Item(16, 2).Value = "concluídas"
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 2
$ws.Cells.Item(17, 1).Value = "Gilvan Pessoa Furtado"
$ws.Cells.Item(17, 2).Value = "concluídas"
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 6
$ws.Cells.Item(18, 1).Value = "Gilvan Pessoa Furtado"
$ws.Cells.Item(18, 2).Value = "em andamento"
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 3
$ws.Cells.Item(19, 1).Value = "Giovanny Augusto Camacho Antevere Mazzarotto"
$ws.Cells.Item(19, 2).Value = "em andamento"
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(20, 1).Value = "Giovanny Augusto Camacho Antevere Mazzarotto"
$ws.Cells.Item(20, 2).Value = "concluídas"
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(21, 1).Value = "Ivana Cristina de Holanda Cunha Barreto"
$ws.Cells.Item(21, 2).Value = "concluídas"
$ws.Cells.Item(21, 3).Value = 11
$ws.Cells.Item(21, 4).Value = 3
$ws.Cells.Item(22, 1).Value = "Ivana Cristina de Holanda Cunha Barreto"
$ws.Cells.Item(22, 2).Value = "em andamento"
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(22, 4).Value = 6
$ws.Cells.Item(23, 1).Value = "Jaime Ribeiro Filho"
$ws.Cells.Item(23, 2).Value = "em andamento"
$ws.Cells.Item(23, 3).Value = 2
$ws.Cells.Item(23, 4).Value = 3
$ws.Cells.Item(24, 1).Value = "Jaime Ribeiro Filho"
$ws.Cells.Item(24, 2).Value = "concluídas"
$ws.Cells.Item(24, 3).Value = 2
$ws.Cells.Item(24, 4).Value = 2
$ws.Cells.Item(25, 1).Value = "João Hermínio Martins da Silva"
$ws.Cells.Item(25, 2).Value = "concluídas"
$ws.Cells.Item(25, 3).Value = 5
$ws.Cells.Item(25, 4).Value = 5
$ws.Cells.Item(26, 1).Value = "João Hermínio Martins da Silva"
$ws.Cells.Item(26, 2).Value = "em andamento"
$ws.Cells.Item(26, 3).Value = 5
$ws.Cells.Item(26, 4).Value = 2
$ws.Cells.Item(27, 1).Value = "Luiz Odorico Monteiro de Andrade"
$ws.Cells.Item(27, 2).Value = "em andamento"
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 4).Value = 5
$ws.Cells.Item(28, 1).Value = "Luiz Odorico Monteiro de Andrade"
$ws.Cells.Item(28, 2).Value = "concluídas"
$ws.Cells.Item(28, 3).Value = 3
$ws.Cells.Item(28, 4).Value = 4
$ws.Cells.Item(29, 1).Value = "Marcela Helena Gambim Fonseca"
$ws.Cells.Item(29, 2).Value = "concluídas"
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(29, 4).Value = 2
$ws.Cells.Item(30, 1).Value = "Marcos Roberto Lourenzoni"
$ws.Cells.Item(30, 2).Value = "em andamento"
$ws.Cells.Item(30, 3).Value = 2
$ws.Cells.Item(30, 4).Value = 9
$ws.Cells.Item(31, 1).Value = "Marcos Roberto Lourenzoni"
$ws.Cells.Item(31, 2).Value = "concluídas"
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(31, 4).Value = 5
$ws.Cells.Item(32, 1).Value = "Maximiliano Ponte"
$ws.Cells.Item(32, 2).Value = "concluídas"
$ws.Cells.Item(32, 3).Value = 1
$ws.Cells.Item(32, 4).Value = 2
$ws.Cells.Item(33, 1).Value = "Maximiliano Ponte"
$ws.Cells.Item(33, 2).Value = "em andamento"
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(34, 1).Value = "Márcio Flávio Moura de Araújo"
$ws.Cells.Item(34, 2).Value = "em andamento"
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = 5
$ws.Cells.Item(35, 1).Value = "Márcio Flávio Moura de Araújo"
$ws.Cells.Item(35, 2).Value = "concluídas"
$ws.Cells.Item(35, 3).Value = 4
$ws.Cells.Item(35, 4).Value = 7
$ws.Cells.Item(36, 1).Value = "Regis Bernardo Brandim Gomes"
$ws.Cells.Item(36, 2).Value = "em andamento"
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(36, 4).Value = 0
$ws.Cells.Item(37, 1).Value = "Regis Bernardo Brandim Gomes"
$ws.Cells.Item(37, 2).Value = "concluídas"
$ws.Cells.Item(37, 3).Value = 6
$ws.Cells.Item(37, 4).Value = 0
$ws.Cells.Item(38, 1).Value = "Roberto Nicolete"
$ws.Cells.Item(38, 2).Value = "concluídas"
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(38, 4).Value = 8
$ws.Cells.Item(39, 1).Value = "Roberto Nicolete"
$ws.Cells.Item(39, 2).Value = "em andamento"
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 4).Value = 5
$ws.Cells.Item(40, 1).Value = "Roberto Wagner Júnior Freire de Freitas"
$ws.Cells.Item(40, 2).Value = "em andamento"
$ws.Cells.Item(40, 3).Value = 5
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(41, 1).Value = "Roberto Wagner Júnior Freire de Freitas"
$ws.Cells.Item(41, 2).Value = "concluídas"
$ws.Cells.Item(41, 3).Value = 7
$ws.Cells.Item(41, 4).Value = 8
$ws.Cells.Item(42, 1).Value = "Sharmênia de Araújo Soares Nuto"
$ws.Cells.Item(42, 2).Value = "concluídas"
$ws.Cells.Item(42, 3).Value = 8
$ws.Cells.Item(42, 4).Value = 3
$ws.Cells.Item(43, 1).Value = "Sharmênia de Araújo Soares Nuto"
$ws.Cells.Item(43, 2).Value = "em andamento"
$ws.Cells.Item(43, 3).Value = 3
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(44, 1).Value = "Vanira Matos Pessoa"
$ws.Cells.Item(44, 2).Value = "concluídas"
$ws.Cells.Item(44, 3).Value = 12
$ws.Cells.Item(44, 4).Value = 1
$ws.Cells.Item(45, 1).Value = "Vanira Matos Pessoa"
$ws.Cells.Item(45, 2).Value = "em andamento"
$ws.Cells.Item(45, 3).Value = 7
$ws.Cells.Item(45, 4).Value = 0
